$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the USD Amount value in T2 (50802 -> 51597)
$ws.Range("T2").Value = 51597

# Move the active selection down to T3 (matches post-edit cursor position)
[void]$ws.Range("T3").Select()
